$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (year 2025) - update recurrence metrics with new totals
$ws.Range("C8").Value = 1287
$ws.Range("E8").Value = 1081
$ws.Range("G8").Value = 83.99378399378399
$ws.Range("H8").Value = 16.00621600621601
